$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4105.647
$ws.Range("I98").Value = 4105.647
$ws.Range("K98").Value = 4105.647
$ws.Range("M98").Value = -2607.647
$ws.Range("H115").Value = 878.8333
$ws.Range("I115").Value = 878.8333
$ws.Range("K115").Value = 2636.4999
$ws.Range("M115").Value = -1069.4999
$ws.Range("H116").Value = 31260374
$ws.Range("I116").Value = 83339000
$ws.Range("K116").Value = 83339000
$ws.Range("M116").Value = -83335558
$ws.Range("H122").Value = 4105.647
$ws.Range("I122").Value = 4105.647
$ws.Range("K122").Value = 12316.941
$ws.Range("M122").Value = -9866.940999999999
$ws.Range("H132").Value = 1051.5518
$ws.Range("I132").Value = 887.9167
$ws.Range("J132").Value = 1837
$ws.Range("K132").Value = 2663.7501
$ws.Range("L132").Value = 5511
$ws.Range("M132").Value = -133.7501000000002
$ws.Range("N132").Value = -10571
$ws.Range("H137").Value = 4214.2856
$ws.Range("I137").Value = 3454.6365
$ws.Range("K137").Value = 10363.9095
$ws.Range("M137").Value = -7813.9095
$ws.Range("H138").Value = 1452767.4
$ws.Range("J138").Value = 2707908.8
$ws.Range("L138").Value = 8123726.399999999
$ws.Range("N138").Value = -8134006.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 40003050
$ws.Range("I61").Value = 1356.5625
$ws.Range("J61").Value = 111117170
$ws.Range("K61").Value = 1356.5625
$ws.Range("L61").Value = 111117170
$ws.Range("M61").Value = -1144.5625
$ws.Range("N61").Value = -111117594
$ws.Range("H74").Value = 25392.28
$ws.Range("I74").Value = 29217.805
$ws.Range("K74").Value = 29217.805
$ws.Range("M74").Value = -28343.805
$ws.Range("H77").Value = 25392.28
$ws.Range("I77").Value = 29217.805
$ws.Range("K77").Value = 146089.025
$ws.Range("M77").Value = -141721.025
$ws.Range("H97").Value = 2875000.2
$ws.Range("I97").Value = 1254.591
$ws.Range("J97").Value = 11906772
$ws.Range("K97").Value = 1254.591
$ws.Range("L97").Value = 11906772
$ws.Range("M97").Value = -758.5909999999999
$ws.Range("N97").Value = -11907764
$ws.Range("H110").Value = 37038010
$ws.Range("I110").Value = 999.8333
$ws.Range("J110").Value = 111112020
$ws.Range("K110").Value = 999.8333
$ws.Range("L110").Value = 111112020
$ws.Range("M110").Value = 1045.1667
$ws.Range("N110").Value = -111116110
$ws.Range("H132").Value = 3583.1858
$ws.Range("I132").Value = 2554.6316
$ws.Range("K132").Value = 7663.8948
$ws.Range("M132").Value = -5133.8948
$ws.Range("H136").Value = 40003050
$ws.Range("I136").Value = 1356.5625
$ws.Range("J136").Value = 111117170
$ws.Range("K136").Value = 4069.6875
$ws.Range("L136").Value = 333351510
$ws.Range("M136").Value = -1519.6875
$ws.Range("N136").Value = -333356610
$ws.Range("H137").Value = 60780
$ws.Range("J137").Value = 60780
$ws.Range("L137").Value = 60780
$ws.Range("N137").Value = -70980

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1619.2559
$ws.Range("I94").Value = 661.6875
$ws.Range("J94").Value = 4404.909
$ws.Range("K94").Value = 661.6875
$ws.Range("L94").Value = 4404.909
$ws.Range("M94").Value = -210.6875
$ws.Range("N94").Value = -5306.909
$ws.Range("H107").Value = 30413746
$ws.Range("I107").Value = 38800616
$ws.Range("K107").Value = 38800616
$ws.Range("M107").Value = -38798696

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4126.6333
$ws.Range("I134").Value = 2223.0527
$ws.Range("J134").Value = 7414.636
$ws.Range("K134").Value = 6669.158100000001
$ws.Range("L134").Value = 22243.908
$ws.Range("M134").Value = -4134.158100000001
$ws.Range("N134").Value = -27313.908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 757.9524
$ws.Range("I5").Value = 513.41174
$ws.Range("J5").Value = 1797.25
$ws.Range("K5").Value = 1540.23522
$ws.Range("L5").Value = 5391.75
$ws.Range("M5").Value = -1428.23522
$ws.Range("N5").Value = -5615.75
$ws.Range("H23").Value = 510.05884
$ws.Range("I23").Value = 349
$ws.Range("J23").Value = 577.1667
$ws.Range("K23").Value = 1047
$ws.Range("L23").Value = 1731.5001
$ws.Range("M23").Value = -812
$ws.Range("N23").Value = -2201.5001
$ws.Range("H135").Value = 757.9524
$ws.Range("I135").Value = 513.41174
$ws.Range("J135").Value = 1797.25
$ws.Range("K135").Value = 4620.70566
$ws.Range("L135").Value = 16175.25
$ws.Range("M135").Value = -2085.70566
$ws.Range("N135").Value = -21245.25
$ws.Range("H141").Value = 3042.923
$ws.Range("I141").Value = 3042.923
$ws.Range("K141").Value = 9128.769
$ws.Range("M141").Value = -3948.769

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 338.6
$ws.Range("I2").Value = 315.5
$ws.Range("J2").Value = 354
$ws.Range("K2").Value = 315.5
$ws.Range("L2").Value = 354
$ws.Range("M2").Value = -202.5
$ws.Range("N2").Value = -580
$ws.Range("H21").Value = 2000
$ws.Range("J21").Value = 2000
$ws.Range("L21").Value = 2000
$ws.Range("N21").Value = -2346
$ws.Range("H29").Value = 20000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 20000
$ws.Range("K29").Value = 0
$ws.Range("L29").ClearContents()
$ws.Range("M29").Value = 20000
$ws.Range("N29").Value = -20580
$ws.Range("H30").Value = 2000
$ws.Range("J30").Value = 2000
$ws.Range("L30").Value = 2000
$ws.Range("N30").Value = -2210
$ws.Range("H31").Value = 2166.6667
$ws.Range("I31").Value = 2166.6667
$ws.Range("K31").Value = 2166.6667
$ws.Range("M31").Value = -1874.6667
$ws.Range("H35").Value = 24999
$ws.Range("J35").Value = 24999
$ws.Range("L35").Value = 24999
$ws.Range("N35").Value = -25595
$ws.Range("H37").Value = 2166.6667
$ws.Range("I37").Value = 2166.6667
$ws.Range("K37").Value = 2166.6667
$ws.Range("M37").Value = -1889.6667
$ws.Range("H40").Value = 10000
$ws.Range("I40").Value = 10000
$ws.Range("K40").Value = 10000
$ws.Range("M40").Value = -9849
$ws.Range("H43").Value = 2017
$ws.Range("I43").Value = 2017
$ws.Range("K43").Value = 2017
$ws.Range("M43").Value = -1866
$ws.Range("H44").Value = 49000
$ws.Range("J44").Value = 49000
$ws.Range("L44").Value = 49000
$ws.Range("N44").Value = -50192
$ws.Range("H48").Value = 11833.333
$ws.Range("I48").Value = 2750
$ws.Range("J48").Value = 30000
$ws.Range("K48").Value = 2750
$ws.Range("L48").Value = 30000
$ws.Range("M48").Value = -2265
$ws.Range("N48").Value = -30970
$ws.Range("H57").Value = 61903.5
$ws.Range("I57").Value = 3750
$ws.Range("K57").Value = 3750
$ws.Range("M57").Value = -2930
$ws.Range("H58").Value = 69752.89999999999
$ws.Range("I58").Value = 3041
$ws.Range("J58").Value = 77165.336
$ws.Range("K58").Value = 3041
$ws.Range("L58").Value = 77165.336
$ws.Range("M58").Value = -2764
$ws.Range("N58").Value = -77719.336
$ws.Range("H102").Value = 3447.342
$ws.Range("I102").Value = 3298.3713
$ws.Range("K102").Value = 3298.3713
$ws.Range("M102").Value = -1676.3713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5894.125
$ws.Range("J7").Value = 6445.8184
$ws.Range("L7").Value = 6445.8184
$ws.Range("N7").Value = -6669.8184
$ws.Range("H55").Value = 320.9
$ws.Range("I55").Value = 54.88889
$ws.Range("K55").Value = 54.88889
$ws.Range("M55").Value = 118.11111
$ws.Range("H122").Value = 3864.0715
$ws.Range("J122").Value = 6531.5
$ws.Range("L122").Value = 19594.5
$ws.Range("N122").Value = -24494.5
$ws.Range("H126").Value = 5894.125
$ws.Range("J126").Value = 6445.8184
$ws.Range("L126").Value = 19337.4552
$ws.Range("N126").Value = -24277.4552
$ws.Range("H132").Value = 7582507.5
$ws.Range("I132").Value = 13891627
$ws.Range("K132").Value = 41674881
$ws.Range("M132").Value = -41672351

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2173.8823
$ws.Range("I132").Value = 996.2308
$ws.Range("J132").Value = 6001.25
$ws.Range("K132").Value = 2988.6924
$ws.Range("L132").Value = 18003.75
$ws.Range("M132").Value = -458.6923999999999
$ws.Range("N132").Value = -23063.75
$ws.Range("H136").Value = 27301330
$ws.Range("I136").Value = 71429940
$ws.Range("K136").Value = 214289820
$ws.Range("M136").Value = -214287270
